# Insert a new data row before the current row 108 ("Hortaliza, Terminal
# Hortofrutícola Agro Chillán - Poroto verde"), shifting the existing rows
# 108-152 down to 109-153 and extending the sheet to R153. Then populate the
# newly inserted row 108 with its final values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 108 - this shifts rows 108:152 down to 109:153
# and pushes the sheet's used range/dimension to A1:R153.
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with the new record's data.
$ws.Cells.Item(108, 1).Value = 7
$ws.Cells.Item(108, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(108, 3).Value = "Ñuble"
$ws.Cells.Item(108, 4).Value = 45119
$ws.Cells.Item(108, 5).Value = 16
$ws.Cells.Item(108, 6).Value = 100112031
$ws.Cells.Item(108, 7).Value = "Poroto verde"
$ws.Cells.Item(108, 8).Value = "Sin especificar"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 20
$ws.Cells.Item(108, 11).Value = 23000
$ws.Cells.Item(108, 12).Value = 23000
$ws.Cells.Item(108, 13).Value = 23000
$ws.Cells.Item(108, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(108, 15).Value = "Perú"
$ws.Cells.Item(108, 16).Value = 920
$ws.Cells.Item(108, 17).Value = 25
$ws.Cells.Item(108, 18).Value = "Hortaliza"
